$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.686.96"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.072.16"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("D5").Value = "'233.23"
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").Value = "'0.622"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D8").Value = "'57.99"
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").Value = "'0.389"
$ws.Range("D10").Value = "'0.0781"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("D12").Value = "2.382.06"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "'14.77"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").Value = "'20.82"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "'0.770"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").Value = "'5.28"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "2.099.53"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "37.660.01"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "'6.18"
$ws.Range("E19").Value = "  -1.00%  "
$ws.Range("D20").Value = "'71.04"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").Value = "0.0₃0830"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "'227.46"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = "  -1.44%  "
$ws.Range("D26").Value = "'169.38"
$ws.Range("E26").Value = "  +0.79%  "
$ws.Range("D27").Value = "'0.138"
$ws.Range("E27").Value = "  +3.41%  "
$ws.Range("D28").Value = "'8.97"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "'19.39"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").Value = "'1.39"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("D32").Value = "'4.66"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "'0.0628"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").Value = "'4.64"
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "'2.45"
$ws.Range("E35").Value = "  -5.32%  "
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("D37").Value = "'3.38"
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").Value = "'5.33"
$ws.Range("E39").Value = "  -5.13%  "
$ws.Range("D40").Value = "'0.0978"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").Value = "'97.83"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  -2.66%  "
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Value = "1.453.05"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").Value = "'16.47"
$ws.Range("E45").Value = "  +5.66%  "
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").Value = "'4.23"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").Value = "'7.36"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").Value = "'3.01"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").Value = "2.266.02"
$ws.Range("E51").Value = "  -1.68%  "
